# Apply updated crypto price/volume figures to the "cryptos" sheet.
# (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.578.44'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '2.518.12'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''313.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = '''98.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.08%  '
$ws.Range('E7').Value = '  -1.32%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '''0.516'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.39%  '
$ws.Range('D10').Value = '''35.15'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('D11').Value = '''0.0799'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('E13').Value = '  -2.70%  '
$ws.Range('D14').Value = '2.902.03'
$ws.Range('E14').Value = '  -1.33%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '''15.21'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.09%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.503.62'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('D17').Value = '''0.807'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.02%  '
$ws.Range('D18').Value = '42.524.90'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').Value = '''6.59'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.97%  '
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').Value = '''12.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('D22').Value = '''68.76'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.42%  '
$ws.Range('D23').Value = '''240.81'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.46%  '
$ws.Range('D24').Value = '''2.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.40%  '
$ws.Range('D25').Value = '''1.99'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.47%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = '''25.32'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.96%  '
$ws.Range('E28').Value = '  -4.43%  '
$ws.Range('D29').Value = '''9.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('E30').Value = '  -6.76%  '
$ws.Range('D31').Value = '''5.85'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.75%  '
$ws.Range('D32').Value = '''156.55'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('E33').Value = '  -2.63%  '
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('D35').Value = '''0.0782'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.53%  '
$ws.Range('D36').Value = '''3.14'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.68%  '
$ws.Range('E37').Value = '  -5.28%  '
$ws.Range('D38').Value = '''17.55'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.96%  '
$ws.Range('E39').Value = '  -3.90%  '
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('D41').Value = '''4.16'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.84%  '
$ws.Range('D42').Value = '''21.77'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('D44').Value = '''0.0295'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('D45').Value = '2.003.76'
$ws.Range('E45').Value = '  +1.71%  '
$ws.Range('D46').Value = '''3.22'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.90%  '
$ws.Range('D47').Value = '''8.92'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('D48').Value = '2.748.92'
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('D49').Value = '''78.71'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.82%  '
$ws.Range('D50').Value = '''0.188'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.51%  '
$ws.Range('D51').Value = '''71.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.63%  '
